$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 14-17 (Resolving-Mac as sender rows duplicated; data recomputed with new TPM values)
$ws.Range("A14:T17").EntireRow.Delete()

# Update Sending cluster (A) and Target cluster (D) labels, and recompute numeric columns G:T

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 95.24255366666667
$ws.Range("H2").Value = 285.727661
$ws.Range("I2").Value = 0.2732032672746878
$ws.Range("J2").Value = 0.2732032672746877
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.597802666666666
$ws.Range("N2").Value = 4.793407999999999
$ws.Range("O2").Value = 0.02304920886321625
$ws.Range("P2").Value = 0.02304920886321625
$ws.Range("Q2").Value = 152.1788062287431
$ws.Range("R2").Value = 1369.609256058688
$ws.Range("S2").Value = 0.006297119169527371
$ws.Range("T2").Value = 0.006297119169527371

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("G3").Value = 95.24255366666667
$ws.Range("H3").Value = 285.727661
$ws.Range("I3").Value = 0.2732032672746878
$ws.Range("J3").Value = 0.2732032672746877
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 1.27306
$ws.Range("N3").Value = 3.81918
$ws.Range("O3").Value = 0.0183646118807784
$ws.Range("P3").Value = 0.0183646118807784
$ws.Range("Q3").Value = 121.2494853708866
$ws.Range("R3").Value = 1091.24536833798
$ws.Range("S3").Value = 0.005017271968060209
$ws.Range("T3").Value = 0.005017271968060209

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 95.24255366666667
$ws.Range("H4").Value = 285.727661
$ws.Range("I4").Value = 0.2732032672746878
$ws.Range("J4").Value = 0.2732032672746877
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 66.45050433333334
$ws.Range("N4").Value = 199.351513
$ws.Range("O4").Value = 0.9585861792560053
$ws.Range("P4").Value = 0.9585861792560054
$ws.Range("Q4").Value = 6328.915725144567
$ws.Range("R4").Value = 56960.2415263011
$ws.Range("S4").Value = 0.2618888761371002
$ws.Range("T4").Value = 0.2618888761371002

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "ECs"
$ws.Range("G5").Value = 113.1680936666667
$ws.Range("H5").Value = 339.504281
$ws.Range("I5").Value = 0.3246226791565123
$ws.Range("J5").Value = 0.3246226791565122
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.597802666666666
$ws.Range("N5").Value = 4.793407999999999
$ws.Range("O5").Value = 0.02304920886321625
$ws.Range("P5").Value = 0.02304920886321625
$ws.Range("Q5").Value = 180.8202818421831
$ws.Range("R5").Value = 1627.382536579648
$ws.Range("S5").Value = 0.007482295933615287
$ws.Range("T5").Value = 0.007482295933615287

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("G6").Value = 113.1680936666667
$ws.Range("H6").Value = 339.504281
$ws.Range("I6").Value = 0.3246226791565123
$ws.Range("J6").Value = 0.3246226791565122
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 1.27306
$ws.Range("N6").Value = 3.81918
$ws.Range("O6").Value = 0.0183646118807784
$ws.Range("P6").Value = 0.0183646118807784
$ws.Range("Q6").Value = 144.0697733232867
$ws.Range("R6").Value = 1296.62795990958
$ws.Range("S6").Value = 0.0059615695104078
$ws.Range("T6").Value = 0.0059615695104078

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("G7").Value = 113.1680936666667
$ws.Range("H7").Value = 339.504281
$ws.Range("I7").Value = 0.3246226791565123
$ws.Range("J7").Value = 0.3246226791565122
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 66.45050433333334
$ws.Range("N7").Value = 199.351513
$ws.Range("O7").Value = 0.9585861792560053
$ws.Range("P7").Value = 0.9585861792560054
$ws.Range("Q7").Value = 7520.076898591907
$ws.Range("R7").Value = 67680.69208732716
$ws.Range("S7").Value = 0.3111788137124892
$ws.Range("T7").Value = 0.3111788137124892

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("D8").Value = "ECs"
$ws.Range("G8").Value = 89.83461266666666
$ws.Range("H8").Value = 269.503838
$ws.Range("I8").Value = 0.2576905883979786
$ws.Range("J8").Value = 0.2576905883979785
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 1.597802666666666
$ws.Range("N8").Value = 4.793407999999999
$ws.Range("O8").Value = 0.02304920886321625
$ws.Range("P8").Value = 0.02304920886321625
$ws.Range("Q8").Value = 143.5379836777671
$ws.Range("R8").Value = 1291.841853099904
$ws.Range("S8").Value = 0.005939564194070098
$ws.Range("T8").Value = 0.005939564194070097

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("D9").Value = "FAPs"
$ws.Range("G9").Value = 89.83461266666666
$ws.Range("H9").Value = 269.503838
$ws.Range("I9").Value = 0.2576905883979786
$ws.Range("J9").Value = 0.2576905883979785
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 1.27306
$ws.Range("N9").Value = 3.81918
$ws.Range("O9").Value = 0.0183646118807784
$ws.Range("P9").Value = 0.0183646118807784
$ws.Range("Q9").Value = 114.3648520014266
$ws.Range("R9").Value = 1029.28366801284
$ws.Range("S9").Value = 0.004732387641258294
$ws.Range("T9").Value = 0.004732387641258294

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("G10").Value = 89.83461266666666
$ws.Range("H10").Value = 269.503838
$ws.Range("I10").Value = 0.2576905883979786
$ws.Range("J10").Value = 0.2576905883979785
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 66.45050433333334
$ws.Range("N10").Value = 199.351513
$ws.Range("O10").Value = 0.9585861792560053
$ws.Range("P10").Value = 0.9585861792560054
$ws.Range("Q10").Value = 5969.555318289655
$ws.Range("R10").Value = 53725.99786460689
$ws.Range("S10").Value = 0.2470186365626502
$ws.Range("T10").Value = 0.2470186365626502

# Row 11
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("D11").Value = "ECs"
$ws.Range("G11").Value = 50.36899566666667
$ws.Range("H11").Value = 151.106987
$ws.Range("I11").Value = 0.1444834651708214
$ws.Range("J11").Value = 0.1444834651708214
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 1.597802666666666
$ws.Range("N11").Value = 4.793407999999999
$ws.Range("O11").Value = 0.02304920886321625
$ws.Range("P11").Value = 0.02304920886321625
$ws.Range("Q11").Value = 80.47971559352177
$ws.Range("R11").Value = 724.3174403416959
$ws.Range("S11").Value = 0.003330229566003494
$ws.Range("T11").Value = 0.003330229566003493

# Row 12
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("D12").Value = "FAPs"
$ws.Range("G12").Value = 50.36899566666667
$ws.Range("H12").Value = 151.106987
$ws.Range("I12").Value = 0.1444834651708214
$ws.Range("J12").Value = 0.1444834651708214
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 1.27306
$ws.Range("N12").Value = 3.81918
$ws.Range("O12").Value = 0.0183646118807784
$ws.Range("P12").Value = 0.0183646118807784
$ws.Range("Q12").Value = 64.12275362340667
$ws.Range("R12").Value = 577.10478261066
$ws.Range("S12").Value = 0.0026533827610521
$ws.Range("T12").Value = 0.0026533827610521

# Row 13
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("G13").Value = 50.36899566666667
$ws.Range("H13").Value = 151.106987
$ws.Range("I13").Value = 0.1444834651708214
$ws.Range("J13").Value = 0.1444834651708214
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 66.45050433333334
$ws.Range("N13").Value = 199.351513
$ws.Range("O13").Value = 0.9585861792560053
$ws.Range("P13").Value = 0.9585861792560054
$ws.Range("Q13").Value = 3347.045164813482
$ws.Range("R13").Value = 30123.40648332133
$ws.Range("S13").Value = 0.1384998528437658
$ws.Range("T13").Value = 0.1384998528437658
